# group key/values added to VfM
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Q1_20_21
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Q1_20_21")

# Clear out the two existing data rows so we can lay the new rows down fresh.
$ws1.Rows.Item(3).ClearContents()
$ws1.Rows.Item(4).ClearContents()

# Row 3 - new "Mars" group summary row
$ws1.Range("B3").Value = "Mars"
$ws1.Range("C3").Value = 28369
$ws1.Range("D3").Value = 14.58
$ws1.Range("E3").Value = 12.98
$ws1.Range("F3").Value = "Very High"
$ws1.Range("G3").Value = "Very High"
$ws1.Range("H3").Value = "Very High"
$ws1.Range("I3").Value = 2089
$ws1.Range("J3").Value = 30458
$ws1.Range("K3").Value = "All you need is love, love is all you need "

# Row 4 - SoT (now tagged with its group in column A)
$ws1.Range("A4").Value = "Roads Places and Environment Group"
$ws1.Range("B4").Value = "SoT"
$ws1.Range("C4").Value = 1469.2
$ws1.Range("D4").Value = 2.58
$ws1.Range("E4").Value = 1.36
$ws1.Range("F4").Value = "High"
$ws1.Range("G4").Value = "High"
$ws1.Range("H4").Value = "High"
$ws1.Range("I4").Value = "928 -678"
$ws1.Range("J4").Value = 2398

# Row 5 - new HSMRPG / A13 row
$ws1.Range("A5").Value = "HSMRPG"
$ws1.Range("B5").Value = "A13"
$ws1.Range("C5").Value = 1985
$ws1.Range("D5").Value = 2.3
$ws1.Range("E5").Value = 2.3
$ws1.Range("F5").Value = "High"
$ws1.Range("I5").Value = 833
$ws1.Range("J5").Value = 3494

# Row 6 - new Rail Group / F9 row
$ws1.Range("A6").Value = "Rail Group"
$ws1.Range("B6").Value = "F9"
$ws1.Range("C6").Value = 1356
$ws1.Range("D6").Value = 1.46
$ws1.Range("E6").Value = 0.74
$ws1.Range("F6").Value = "Medium"
$ws1.Range("G6").Value = "N/A"
$ws1.Range("H6").Value = "N/A"
$ws1.Range("I6").Value = 2956
$ws1.Range("J6").Value = 4312

# Row 7 - Columbia (now tagged with its group in column A)
$ws1.Range("A7").Value = "Roads Places and Environment Group"
$ws1.Range("B7").Value = "Columbia"
$ws1.Range("D7").Value = 0.38
$ws1.Range("E7").Value = 0.63
$ws1.Range("F7").Value = "Poor"
$ws1.Range("I7").Value = 1172
$ws1.Range("J7").Value = 738.36

# ---------------------------------------------------------------------------
# Sheet 2: Q4_19_20
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Q4_19_20")

$ws2.Rows.Item(3).ClearContents()
$ws2.Rows.Item(4).ClearContents()

# Row 3 - new "Mars" group summary row
$ws2.Range("B3").Value = "Mars"
$ws2.Range("C3").Value = 30292.2
$ws2.Range("D3").Value = 21.45
$ws2.Range("E3").Value = 19.72
$ws2.Range("F3").Value = "Very High"
$ws2.Range("I3").Value = 1481.6
$ws2.Range("J3").Value = 31773.8
$ws2.Range("K3").Value = "Hello is it me you’re looking for"

# Row 4 - SoT (now tagged with its group in column A)
$ws2.Range("A4").Value = "Roads Places and Environment Group"
$ws2.Range("B4").Value = "SoT"
$ws2.Range("C4").Value = 1469.2
$ws2.Range("D4").Value = 2.58
$ws2.Range("E4").Value = 1.36
$ws2.Range("F4").Value = "High"
$ws2.Range("I4").Value = 928
$ws2.Range("J4").Value = 2398
$ws2.Range("K4").Value = "Please allow me to introduce myself I’m a man of wealth and taste."

# Row 5 - new HSMRPG / A11 row (only group + key populated)
$ws2.Range("A5").Value = "HSMRPG"
$ws2.Range("B5").Value = "A11"

# Row 6 - new Rail Group / A13 row
$ws2.Range("A6").Value = "Rail Group"
$ws2.Range("B6").Value = "A13"
$ws2.Range("C6").Value = 1985
$ws2.Range("D6").Value = 2.3
$ws2.Range("E6").Value = 2.3
$ws2.Range("F6").Value = "High"
$ws2.Range("I6").Value = 833
$ws2.Range("J6").Value = 3494

# Row 7 - F9 (now tagged with its group in column A)
$ws2.Range("A7").Value = "Roads Places and Environment Group"
$ws2.Range("B7").Value = "F9"
$ws2.Range("C7").Value = 2952
$ws2.Range("D7").Value = 1.54
$ws2.Range("E7").Value = 0.78
$ws2.Range("F7").Value = "Medium"
$ws2.Range("I7").Value = 2831
$ws2.Range("J7").Value = 4364

# Row 8 - new Rail Group / Columbia row
$ws2.Range("A8").Value = "Rail Group"
$ws2.Range("B8").Value = "Columbia"
$ws2.Range("D8").Value = 0.38
$ws2.Range("E8").Value = 0.63
$ws2.Range("F8").Value = "Poor"
$ws2.Range("I8").Value = 1172
$ws2.Range("J8").Value = 738.36

# ---------------------------------------------------------------------------
# Sheet 3: Count - recomputed per-category totals / counts now that the
# group rows feed into the PVC totals and category counts.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Count")

# PVC total per category (Q1 20/21 = column C, Q4 19/20 = column D)
$ws3.Range("D4").Value = 1172      # Poor
$ws3.Range("C6").Value = 2956      # Medium
$ws3.Range("C7").Value = 833       # High
$ws3.Range("D7").Value = 1761      # High
$ws3.Range("C8").Value = 2089      # Very High
$ws3.Range("D8").Value = 1481.6    # Very High
$ws3.Range("C11").Value = 7050     # Total
$ws3.Range("D11").Value = 7245.6   # Total

# Category count (Q1 20/21 = column C, Q4 19/20 = column D)
$ws3.Range("D16").Value = 1        # Poor
$ws3.Range("C18").Value = 1        # Medium
$ws3.Range("C19").Value = 2        # High
$ws3.Range("D19").Value = 2        # High
$ws3.Range("C20").Value = 1        # Very High
$ws3.Range("D20").Value = 1        # Very High
$ws3.Range("C23").Value = 5        # Total
$ws3.Range("D23").Value = 5        # Total
